$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A244").Value = "IMX-USD"
$ws.Range("A245").Value = "MNT-USD"
$ws.Range("A246").Value = "GRT-USD"
